# ProgramList.xlsx - "Implemented retry logic for all failed test cases"
#
# The "Program" worksheet holds an alphabetically-sorted list of program
# names in column A (column B only carries the "Tech Check Link" header).
# This change keeps that ordering while:
#   1. Adding a new program "Demo 1"      (sorts just before "Dummy Prog 18/5")
#   2. Adding a new program "Importabc"   (sorts just before "Imported 1")
#   3. Adding a new program "LPrgm2"      (sorts just before "LProgram1")
#   4. Removing the obsolete program "JulyUATTest"
#
# Net effect: the used range grows from A1:B44 to A1:B46 (44 -> 46 data rows
# incl. header), i.e. +3 new rows, -1 removed row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program")

# --- 1. Insert "Demo 1" right before "Dummy Prog 18/5" (row 20) ------------
$ws.Rows.Item(20).Insert()
$ws.Range("A20").Value = "Demo 1"

# --- 2. Insert "Importabc" right before "Imported 1" (now row 25) ----------
$ws.Rows.Item(25).Insert()
$ws.Range("A25").Value = "Importabc"

# --- 3. Remove the obsolete "JulyUATTest" row (now row 28) -----------------
$ws.Rows.Item(28).Delete()

# --- 4. Insert "LPrgm2" right before "LProgram1" (now row 29) --------------
$ws.Rows.Item(29).Insert()
$ws.Range("A29").Value = "LPrgm2"
